# "add nipah virus result"
#
# This edit:
#  1. Collapses the detailed Nipah gene list in column V ("Genes (PM)")
#     down to the generic value "Other" for a set of rows.
#  2. Inserts the "Surveillance of Nipah virus in Pteropus medius..."
#     paper ahead of the "Clinico-epidemiological presentations..."
#     paper, which (because both rows already existed) manifests as the
#     two rows (25 and 26) swapping their paper-specific column values.
#  3. Records "NA" for the previously-blank Genes (PM) value on row 27.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- helper: swap the contents of two single cells, preserving the
#     original value "type" (numeric-looking strings such as PMIDs stay
#     text) and correctly clearing a destination whose source was blank.
function Swap-CellValues($ws, $addr1, $addr2, $scratch) {
    $v1 = $ws.Range($addr1).Value()
    $v2 = $ws.Range($addr2).Value()

    if (-not [string]::IsNullOrEmpty($v1)) {
        $ws.Range($addr1).Copy()
        $ws.Range($scratch).PasteSpecial(-4163)
    }

    if ([string]::IsNullOrEmpty($v2)) {
        $ws.Range($addr1).ClearContents()
    } else {
        $ws.Range($addr2).Copy()
        $ws.Range($addr1).PasteSpecial(-4163)
    }

    if ([string]::IsNullOrEmpty($v1)) {
        $ws.Range($addr2).ClearContents()
    } else {
        $ws.Range($scratch).Copy()
        $ws.Range($addr2).PasteSpecial(-4163)
        $ws.Range($scratch).ClearContents()
    }
}

# 1. Genes (PM) -> "Other"
$otherRows = @(5, 7, 8, 9, 10, 12, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24, 29, 30, 31, 32, 34)
foreach ($r in $otherRows) {
    $ws.Range("V$r").Value = "Other"
}

# 2. Swap the paper-specific fields between row 25 and row 26.
$scratch = "ZZ1000"
$swapCols = @("A", "B", "C", "E", "L", "N", "P", "X", "Z", "AB")
foreach ($col in $swapCols) {
    Swap-CellValues $ws "$col`25" "$col`26" $scratch
}

# 3. Genes (PM) for row 27 was blank, now "NA".
$ws.Range("V27").Value = "NA"
